$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("A4").Value = 131009286
$ws.Range("B4").Value = 80252
$ws.Range("D4").Value = 'LC'
$ws.Range("E4").Value = 6456
$ws.Range("F4").Value = 'Skinnlav'
$ws.Range("G4").Value = 'Leptogium saturninum'
$ws.Range("H4").Value = '(Dicks.) Nyl.'
$ws.Range("M4").Value = ""
$ws.Range("Q4").Value = 589865
$ws.Range("R4").Value = 6911173
$ws.Range("Z4").Value = '10:53'
$ws.Range("AB4").Value = '10:53'
$ws.Range("AC4").Value = ""

# Row 5
$ws.Range("A5").Value = 131009298
$ws.Range("Q5").Value = 589735
$ws.Range("R5").Value = 6911227
$ws.Range("Z5").Value = '10:12'
$ws.Range("AB5").Value = '10:12'

# Row 6
$ws.Range("A6").Value = 131009294
$ws.Range("B6").Value = 57884
$ws.Range("D6").Value = 'NT'
$ws.Range("E6").Value = 100109
$ws.Range("F6").Value = 'Tretåig hackspett'
$ws.Range("G6").Value = 'Picoides tridactylus'
$ws.Range("H6").Value = '(Linnaeus, 1758)'
$ws.Range("M6").Value = 'färska spår'
$ws.Range("Q6").Value = 589787
$ws.Range("R6").Value = 6911183
$ws.Range("Z6").Value = '10:24'
$ws.Range("AB6").Value = '10:24'
$ws.Range("AC6").Value = 'färska ringhack på tall'

# Row 7
$ws.Range("A7").Value = 131009283
$ws.Range("B7").Value = 79243
$ws.Range("E7").Value = 6425
$ws.Range("F7").Value = 'Garnlav'
$ws.Range("G7").Value = 'Alectoria sarmentosa'
$ws.Range("H7").Value = '(Ach.) Ach.'
$ws.Range("M7").Value = ""
$ws.Range("Q7").Value = 589968
$ws.Range("R7").Value = 6911120
$ws.Range("Z7").Value = '11:03'
$ws.Range("AB7").Value = '11:03'
$ws.Range("AC7").Value = ""

# Row 8
$ws.Range("A8").Value = 131009297
$ws.Range("Q8").Value = 589752
$ws.Range("R8").Value = 6911214
$ws.Range("Z8").Value = '10:18'
$ws.Range("AB8").Value = '10:18'

# Row 9
$ws.Range("A9").Value = 131009287
$ws.Range("Q9").Value = 589835
$ws.Range("R9").Value = 6911210
$ws.Range("Z9").Value = '10:47'
$ws.Range("AB9").Value = '10:47'

# Row 10
$ws.Range("A10").Value = 131009270
$ws.Range("B10").Value = 57884
$ws.Range("E10").Value = 100109
$ws.Range("F10").Value = 'Tretåig hackspett'
$ws.Range("G10").Value = 'Picoides tridactylus'
$ws.Range("H10").Value = '(Linnaeus, 1758)'
$ws.Range("M10").Value = 'färska spår'
$ws.Range("Q10").Value = 589700
$ws.Range("R10").Value = 6911274
$ws.Range("Z10").Value = '12:22'
$ws.Range("AB10").Value = '12:22'
$ws.Range("AC10").Value = 'färska ringhack på tall'

# Row 21
$ws.Range("A21").Value = 131009308
$ws.Range("B21").Value = 79243
$ws.Range("E21").Value = 6425
$ws.Range("F21").Value = 'Garnlav'
$ws.Range("G21").Value = 'Alectoria sarmentosa'
$ws.Range("H21").Value = '(Ach.) Ach.'
$ws.Range("M21").Value = ""
$ws.Range("Q21").Value = 589686
$ws.Range("R21").Value = 6911077
$ws.Range("Z21").Value = '09:17'
$ws.Range("AB21").Value = '09:17'
$ws.Range("AC21").Value = ""

# Row 22
$ws.Range("A22").Value = 131009305
$ws.Range("B22").Value = 57884
$ws.Range("E22").Value = 100109
$ws.Range("F22").Value = 'Tretåig hackspett'
$ws.Range("G22").Value = 'Picoides tridactylus'
$ws.Range("H22").Value = '(Linnaeus, 1758)'
$ws.Range("M22").Value = 'färska spår'
$ws.Range("Q22").Value = 589757
$ws.Range("R22").Value = 6911178
$ws.Range("Z22").Value = '09:43'
$ws.Range("AB22").Value = '09:43'
$ws.Range("AC22").Value = 'färska ringhack på tall'

# Row 23
$ws.Range("A23").Value = 131009498
$ws.Range("B23").Value = 57884
$ws.Range("E23").Value = 100109
$ws.Range("F23").Value = 'Tretåig hackspett'
$ws.Range("G23").Value = 'Picoides tridactylus'
$ws.Range("H23").Value = '(Linnaeus, 1758)'
$ws.Range("I23").Value = ""
$ws.Range("M23").Value = 'färska spår'
$ws.Range("Q23").Value = 589629
$ws.Range("R23").Value = 6911040
$ws.Range("Z23").Value = '09:02'
$ws.Range("AB23").Value = '09:02'
$ws.Range("AC23").Value = 'färska ringhack på tall'

# Row 24
$ws.Range("A24").Value = 131009304
$ws.Range("B24").Value = 58043
$ws.Range("E24").Value = 103021
$ws.Range("F24").Value = 'Talltita'
$ws.Range("G24").Value = 'Poecile montanus'
$ws.Range("H24").Value = '(Conrad von Baldenstein, 1827)'
$ws.Range("I24").NumberFormat = "@"
$ws.Range("I24").Value = '1'
$ws.Range("M24").Value = 'lockläte, övriga läten'
$ws.Range("Q24").Value = 589753
$ws.Range("R24").Value = 6911167
$ws.Range("Z24").Value = '09:45'
$ws.Range("AB24").Value = '09:45'
$ws.Range("AC24").Value = ""

# Row 39
$ws.Range("A39").Value = 131009291
$ws.Range("B39").Value = 80252
$ws.Range("D39").Value = 'LC'
$ws.Range("E39").Value = 6456
$ws.Range("F39").Value = 'Skinnlav'
$ws.Range("G39").Value = 'Leptogium saturninum'
$ws.Range("H39").Value = '(Dicks.) Nyl.'
$ws.Range("M39").Value = ""
$ws.Range("Q39").Value = 589791
$ws.Range("R39").Value = 6911200
$ws.Range("Z39").Value = '10:37'
$ws.Range("AB39").Value = '10:37'
$ws.Range("AC39").Value = ""

# Row 40
$ws.Range("A40").Value = 131009275
$ws.Range("B40").Value = 57884
$ws.Range("D40").Value = 'NT'
$ws.Range("E40").Value = 100109
$ws.Range("F40").Value = 'Tretåig hackspett'
$ws.Range("G40").Value = 'Picoides tridactylus'
$ws.Range("H40").Value = '(Linnaeus, 1758)'
$ws.Range("M40").Value = 'färska spår'
$ws.Range("Q40").Value = 589844
$ws.Range("R40").Value = 6911365
$ws.Range("Z40").Value = '11:53'
$ws.Range("AB40").Value = '11:53'
$ws.Range("AC40").Value = 'färska ringhack på tall'
